$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily-push row goes in right after the last "2026/02/04" block (row 771),
# which pushes every subsequent row (772-813) down by one to (773-814).
$ws.Rows("772:772").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Force the date column to be stored as literal text (matching every other
# row in the sheet, e.g. A771 = "2026/02/04" as text, not a real date), then
# reset the style back to Normal so no stray number-format style lingers.
$ws.Range("A772").NumberFormat = "@"
$ws.Range("A772").Value = "2026/02/04"
$ws.Range("A772").Style = "Normal"

$ws.Range("B772").Value = "水"
$ws.Range("C772").Value = 23
$ws.Range("D772").Value = 201
